$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- SMOTE update -> SMOTE + Threshold finalized ---
# Update model order/labels (KNN <-> Random Forest swapped within each
# result table) and refresh all metric values across the result blocks.
$ws.Range("B3").Value = 0.95599999999999996
$ws.Range("C3").Value = 0.89800000000000002
$ws.Range("D3").Value = 0.88900000000000001
$ws.Range("E3").Value = 0.876
$ws.Range("F3").Value = 0.86299999999999999
$ws.Range("A4").Value = "Random Forest"
$ws.Range("B4").Value = 0.95199999999999996
$ws.Range("C4").Value = 0.89400000000000002
$ws.Range("D4").Value = 0.82199999999999995
$ws.Range("E4").Value = 0.84199999999999997
$ws.Range("F4").Value = 0.82499999999999996
$ws.Range("A5").Value = "KNN"
$ws.Range("B5").Value = 0.89900000000000002
$ws.Range("C5").Value = 0.76300000000000001
$ws.Range("D5").Value = 0.55600000000000005
$ws.Range("E5").Value = 0.61899999999999999
$ws.Range("F5").Value = 0.58599999999999997
$ws.Range("B6").Value = 0.93400000000000005
$ws.Range("C6").Value = 0.82
$ws.Range("D6").Value = 0.8
$ws.Range("E6").Value = 0.79200000000000004
$ws.Range("F6").Value = 0.76800000000000002
$ws.Range("B7").Value = 0.93
$ws.Range("C7").Value = 0.80200000000000005
$ws.Range("D7").Value = 0.73299999999999998
$ws.Range("E7").Value = 0.76200000000000001
$ws.Range("F7").Value = 0.72499999999999998
$ws.Range("A10").Value = "Regular model -Balance"
$ws.Range("B12").Value = 0.88900000000000001
$ws.Range("C12").Value = 0.93799999999999994
$ws.Range("D12").Value = 0.86699999999999999
$ws.Range("E12").Value = 0.88300000000000001
$ws.Range("F12").Value = 0.80500000000000005
$ws.Range("A13").Value = "Random Forest"
$ws.Range("B13").Value = 0.88900000000000001
$ws.Range("C13").Value = 0.88500000000000001
$ws.Range("D13").Value = 0.91100000000000003
$ws.Range("E13").Value = 0.89300000000000002
$ws.Range("F13").Value = 0.78800000000000003
$ws.Range("A14").Value = "KNN"
$ws.Range("B14").Value = 0.74399999999999999
$ws.Range("C14").Value = 0.872
$ws.Range("D14").Value = 0.6
$ws.Range("E14").Value = 0.69
$ws.Range("F14").Value = 0.52500000000000002
$ws.Range("A16").Value = "XGBoost"
$ws.Range("B16").Value = 0.755
$ws.Range("C16").Value = 0.73699999999999999
$ws.Range("D16").Value = 0.82199999999999995
$ws.Range("E16").Value = 0.71099999999999997
$ws.Range("F16").Value = 0.52600000000000002
$ws.Range("A23").Value = "Regular model -Imbalance"
$ws.Range("A24").Value = "Decision Tree"
$ws.Range("B24").Value = 0.96199999999999997
$ws.Range("C24").Value = 0.91300000000000003
$ws.Range("D24").Value = 0.86699999999999999
$ws.Range("E24").Value = 0.877
$ws.Range("F24").Value = 0.86399999999999999
$ws.Range("A25").Value = "Random Forest"
$ws.Range("B25").Value = 0.96099999999999997
$ws.Range("C25").Value = 0.97499999999999998
$ws.Range("D25").Value = 0.77800000000000002
$ws.Range("E25").Value = 0.86299999999999999
$ws.Range("F25").Value = 0.79500000000000004
$ws.Range("A26").Value = "KNN"
$ws.Range("B26").Value = 0.89200000000000002
$ws.Range("C26").Value = 0.92
$ws.Range("D26").Value = 0.35599999999999998
$ws.Range("E26").Value = 0.497
$ws.Range("F26").Value = 0.52200000000000002
$ws.Range("A27").Value = "SVM"
$ws.Range("A28").Value = " XGBoost"
$ws.Range("B28").Value = 0.89900000000000002
$ws.Range("C28").Value = 0.84699999999999998
$ws.Range("D28").Value = 0.46600000000000003
$ws.Range("E28").Value = 0.58199999999999996
$ws.Range("F28").Value = 0.57199999999999995
$ws.Range("A32").Value = "Threshold only"
$ws.Range("B34").Value = 0.93100000000000005
$ws.Range("C34").Value = 0.73399999999999999
$ws.Range("D34").Value = 0.88900000000000001
$ws.Range("E34").Value = 0.80200000000000005
$ws.Range("F34").Value = 0.76700000000000002
$ws.Range("A35").Value = "Random Forest"
$ws.Range("B35").Value = 0.96499999999999997
$ws.Range("C35").Value = 0.95199999999999996
$ws.Range("D35").Value = 0.82199999999999995
$ws.Range("E35").Value = 0.88
$ws.Range("F35").Value = 0.86399999999999999
$ws.Range("A36").Value = "KNN"
$ws.Range("B36").Value = 0.88900000000000001
$ws.Range("C36").Value = 0.68600000000000005
$ws.Range("D36").Value = 0.57799999999999996
$ws.Range("E36").Value = 0.61299999999999999
$ws.Range("F36").Value = 0.51600000000000001
$ws.Range("A38").Value = "XGBoost"
$ws.Range("A41").Value = "SMOTE only"
$ws.Range("A43").Value = "Decision tree"
$ws.Range("B43").Value = 0.875
$ws.Range("C43").Value = 0.61099999999999999
$ws.Range("D43").Value = 0.6
$ws.Range("E43").Value = 0.58299999999999996
$ws.Range("F43").Value = 0.52500000000000002
$ws.Range("A44").Value = "Random Forest"
$ws.Range("B44").Value = 0.95099999999999996
$ws.Range("C44").Value = 0.89400000000000002
$ws.Range("D44").Value = 0.82199999999999995
$ws.Range("E44").Value = 0.84199999999999997
$ws.Range("F44").Value = 0.82399999999999995
$ws.Range("A45").Value = "KNN"
$ws.Range("B45").Value = 0.875
$ws.Range("C45").Value = 0.61099999999999999
$ws.Range("D45").Value = 0.6
$ws.Range("E45").Value = 0.58299999999999996
$ws.Range("F45").Value = 0.52500000000000002
$ws.Range("D47").Value = 0.71099999999999997
$ws.Range("E47").Value = 0.755
$ws.Range("F47").Value = 0.72499999999999998
$ws.Range("A53").Value = "Decision tree"

# Update the active cell selection to match the final saved state
$ws.Range("K11").Select() | Out-Null
